# Update the "lastlogin" timestamp for the "vikrant" user row (row 4, column E)
# from "2022-01-10 00:36:18.36S" to "2022-01-10 15:38:37.38S".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E4").Value = "2022-01-10 15:38:37.38S"
